$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New filament/printer data for rows 7, 9 and 10 -----------------------
# NOTE: the order these are typed in matters, because it determines the
# order new entries are appended to the shared-strings table.
# Required shared-string append order: "MK3 Pretty PETG V2 rigid.ink",
# "MK3 Pretty PETG V2 PrimaSelect", "EasyWood PLA".
$ws.Range("C10").Value = "MK3 Pretty PETG V2 rigid.ink"
$ws.Range("D10").Value = "MK3 Pretty PETG V2"
$ws.Range("B10").Value = "MK3 Pretty PETG V2"
$ws.Range("B10").ClearFormats()

$ws.Range("C9").Value = "MK3 Pretty PETG V2 PrimaSelect"
$ws.Range("D9").Value = "MK3 Pretty PETG V2"
$ws.Range("B9").Value = "MK3 Pretty PETG V2"
$ws.Range("B9").ClearFormats()

$ws.Range("C7").Value = "EasyWood PLA"
$ws.Range("D7").Value = "EasyWood PLA"

# --- Column widths: split the old A:C block so C gets its own width -------
$ws.Columns("C:C").ColumnWidth = 29.14

# --- Selection shown when the sheet was last saved -------------------------
[void]$ws.Range("G4").Select()
